$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.266.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.914.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.13%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4628'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.77%  '
$ws.Range('E8').Value = '  +2.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.76'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07943'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.36'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.903.54'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.118'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.94%  '
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06953'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.281.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.358'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.10'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.154.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.062'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.94'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.148'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.995'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '118.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09393'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9260'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.360'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.356'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.274'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.200'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05840'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.973'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5753'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.07%  '
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.970'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.305'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5424'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.07074'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.73%  '
$ws.Range('E49').Value = '  +3.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.571'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '113.54'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.04%  '
